# Update cryptos list with new price/volume data (per commit: Tue Apr  4 14:28:10 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.177.50"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.867.05"
$ws.Range("E3").Value = "  +3.12%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'311.80"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5000"
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("D8").Value = "'0.3938"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.09901"
$ws.Range("E9").Value = "  +26.60%  "
$ws.Range("D11").Value = "'41.17"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "'6.459"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "'20.92"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").Value = "1.867.72"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "'1.003"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'7.392"
$ws.Range("D17").Value = "'0.00001135"
$ws.Range("E17").Value = "  +5.21%  "
$ws.Range("D18").Value = "'93.51"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'0.06640"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'17.41"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "'6.108"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "28.248.14"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "'11.33"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "'2.270"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "'2.551"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").Value = "'21.24"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("D28").Value = "2.081.71"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").Value = "'158.02"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").Value = "'128.15"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("D32").Value = "'1.054"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").Value = "'5.614"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "'3.613"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "'0.06804"
$ws.Range("E35").Value = "  -4.69%  "
$ws.Range("D36").Value = "'9.426"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").Value = "'0.02399"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").Value = "'0.2182"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "'5.004"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'11.45"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").Value = "'0.6287"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").Value = "'1.174"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'13.49"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").Value = "'0.6002"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "'3.664"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'1.273"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "'124.65"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "'1.985"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").Value = "'1.199"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "'1.122"
$ws.Range("E51").Value = "  +5.29%  "
